$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8154.3125
$ws.Range("I40").Value = 18216.5
$ws.Range("J40").Value = 2117
$ws.Range("K40").Value = 18216.5
$ws.Range("L40").Value = 2117
$ws.Range("M40").Value = -18041.5
$ws.Range("N40").Value = -2467
$ws.Range("H48").Value = 1450
$ws.Range("J48").Value = 3250
$ws.Range("L48").Value = 9750
$ws.Range("N48").Value = -10334
$ws.Range("H56").Value = 1450
$ws.Range("J56").Value = 3250
$ws.Range("L56").Value = 9750
$ws.Range("N56").Value = -10818
$ws.Range("H106").Value = 260625.67
$ws.Range("I106").Value = 2050
$ws.Range("J106").Value = 777777
$ws.Range("K106").Value = 2050
$ws.Range("L106").Value = 777777
$ws.Range("M106").Value = -1419
$ws.Range("N106").Value = -779039
$ws.Range("H125").Value = 1723.9474
$ws.Range("I125").Value = 1824
$ws.Range("J125").Value = 1697.2667
$ws.Range("K125").Value = 16416
$ws.Range("L125").Value = 15275.4003
$ws.Range("M125").Value = -13956
$ws.Range("N125").Value = -20195.4003
$ws.Range("H138").Value = 1474.89
$ws.Range("I138").Value = 814.8182
$ws.Range("K138").Value = 2444.4546
$ws.Range("M138").Value = 2695.5454

# ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5019.1836
$ws.Range("I32").Value = 2923.7896
$ws.Range("J32").Value = 12257.818
$ws.Range("K32").Value = 2923.7896
$ws.Range("L32").Value = 12257.818
$ws.Range("M32").Value = -2636.7896
$ws.Range("N32").Value = -12831.818
$ws.Range("H46").Value = 83345290
$ws.Range("I46").Value = 500004500
$ws.Range("J46").Value = 13444.6
$ws.Range("K46").Value = 500004500
$ws.Range("L46").Value = 13444.6
$ws.Range("M46").Value = -500004181
$ws.Range("N46").Value = -14082.6
$ws.Range("H63").Value = 3688.9285
$ws.Range("I63").Value = 2408.125
$ws.Range("J63").Value = 5396.6665
$ws.Range("K63").Value = 2408.125
$ws.Range("L63").Value = 5396.6665
$ws.Range("M63").Value = -1722.125
$ws.Range("N63").Value = -6768.6665
$ws.Range("H66").Value = 3688.9285
$ws.Range("I66").Value = 2408.125
$ws.Range("J66").Value = 5396.6665
$ws.Range("K66").Value = 12040.625
$ws.Range("L66").Value = 26983.3325
$ws.Range("M66").Value = -8608.625
$ws.Range("N66").Value = -33847.3325
$ws.Range("H97").Value = 30303800
$ws.Range("I97").Value = 33334016
$ws.Range("J97").Value = 1633.6666
$ws.Range("K97").Value = 33334016
$ws.Range("L97").Value = 1633.6666
$ws.Range("M97").Value = -33333520
$ws.Range("N97").Value = -2625.6666
$ws.Range("H122").Value = 1770.8334
$ws.Range("I122").Value = 1527.8235
$ws.Range("J122").Value = 2361
$ws.Range("K122").Value = 4583.470499999999
$ws.Range("L122").Value = 7083
$ws.Range("M122").Value = -2133.470499999999
$ws.Range("N122").Value = -11983

# BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2399.8462
$ws.Range("I20").Value = 1739.5
$ws.Range("J20").Value = 2693.3333
$ws.Range("K20").Value = 1739.5
$ws.Range("L20").Value = 2693.3333
$ws.Range("M20").Value = -1492.5
$ws.Range("N20").Value = -3187.3333
$ws.Range("H22").Value = 200.5
$ws.Range("I22").Value = 200.5
$ws.Range("K22").Value = 200.5
$ws.Range("M22").Value = -27.5
$ws.Range("H86").Value = 3156.6365
$ws.Range("I86").Value = 2858.2222
$ws.Range("J86").Value = 4499.5
$ws.Range("K86").Value = 2858.2222
$ws.Range("L86").Value = 4499.5
$ws.Range("M86").Value = -1735.2222
$ws.Range("N86").Value = -6745.5
$ws.Range("H89").Value = 3156.6365
$ws.Range("I89").Value = 2858.2222
$ws.Range("J89").Value = 4499.5
$ws.Range("K89").Value = 14291.111
$ws.Range("L89").Value = 22497.5
$ws.Range("M89").Value = -8675.111000000001
$ws.Range("N89").Value = -33729.5
$ws.Range("H94").Value = 1073.081
$ws.Range("I94").Value = 1051.7693
$ws.Range("J94").Value = 1123.4546
$ws.Range("K94").Value = 1051.7693
$ws.Range("L94").Value = 1123.4546
$ws.Range("M94").Value = -600.7692999999999
$ws.Range("N94").Value = -2025.4546
$ws.Range("H107").Value = 1829.4615
$ws.Range("I107").Value = 1433.875
$ws.Range("J107").Value = 2462.4
$ws.Range("K107").Value = 1433.875
$ws.Range("L107").Value = 2462.4
$ws.Range("M107").Value = 486.125
$ws.Range("N107").Value = -6302.4
$ws.Range("H134").Value = 2953.1042
$ws.Range("I134").Value = 1974.0834
$ws.Range("J134").Value = 3279.4443
$ws.Range("K134").Value = 5922.2502
$ws.Range("L134").Value = 9838.332900000001
$ws.Range("M134").Value = -3387.2502
$ws.Range("N134").Value = -14908.3329

# CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9015515
$ws.Range("I31").Value = 2847.4119
$ws.Range("J31").Value = 16676282
$ws.Range("K31").Value = 2847.4119
$ws.Range("L31").Value = 16676282
$ws.Range("M31").Value = -2552.4119
$ws.Range("N31").Value = -16676872
$ws.Range("H34").Value = 9015515
$ws.Range("I34").Value = 2847.4119
$ws.Range("J34").Value = 16676282
$ws.Range("K34").Value = 2847.4119
$ws.Range("L34").Value = 16676282
$ws.Range("M34").Value = -2645.4119
$ws.Range("N34").Value = -16676686
$ws.Range("H105").Value = 3125.2666
$ws.Range("I105").Value = 3687.9
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 3687.9
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -1940.9
$ws.Range("N105").Value = -5494

# CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3838.3
$ws.Range("I68").Value = 876
$ws.Range("K68").Value = 2628
$ws.Range("M68").Value = -1817
$ws.Range("H70").Value = 5144.5713
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630
$ws.Range("H71").Value = 3838.3
$ws.Range("I71").Value = 876
$ws.Range("K71").Value = 7884
$ws.Range("M71").Value = -3828
$ws.Range("H73").Value = 5144.5713
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184
$ws.Range("H107").Value = 12100
$ws.Range("I107").Value = 10394.9
$ws.Range("J107").Value = 14535.857
$ws.Range("K107").Value = 31184.7
$ws.Range("L107").Value = 43607.571
$ws.Range("M107").Value = -29264.7
$ws.Range("N107").Value = -47447.571

# GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5268.421
$ws.Range("I70").Value = 5270.5884
$ws.Range("J70").Value = 5250
$ws.Range("K70").Value = 5270.5884
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = -5000.5884
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 5268.421
$ws.Range("I73").Value = 5270.5884
$ws.Range("J73").Value = 5250
$ws.Range("K73").Value = 5270.5884
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = -4334.5884
$ws.Range("N73").Value = -7122

# LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2573
$ws.Range("I7").Value = 2083.6667
$ws.Range("J7").Value = 3551.6667
$ws.Range("K7").Value = 2083.6667
$ws.Range("L7").Value = 3551.6667
$ws.Range("M7").Value = -1971.6667
$ws.Range("N7").Value = -3775.6667
$ws.Range("H126").Value = 2573
$ws.Range("I126").Value = 2083.6667
$ws.Range("J126").Value = 3551.6667
$ws.Range("K126").Value = 6251.000100000001
$ws.Range("L126").Value = 10655.0001
$ws.Range("M126").Value = -3781.000100000001
$ws.Range("N126").Value = -15595.0001

# WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 865267.7
$ws.Range("I136").Value = 1111954.4
$ws.Range("K136").Value = 3335863.2
$ws.Range("M136").Value = -3333313.2
